$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("566:568").Insert()

$ws.Range("A566").Value = 11
$ws.Range("B566").Value = "Vega Monumental Concepción"
$ws.Range("C566").Value = "Bíobío"
$ws.Range("D566").Value = 44964
$ws.Range("E566").Value = 8
$ws.Range("F566").Value = "Fruta"
$ws.Range("G566").Value = 100108
$ws.Range("H566").Value = "Tropicales y subtropicales"
$ws.Range("I566").Value = 100108006
$ws.Range("J566").Value = "Plátano"
$ws.Range("K566").Value = "Sin especificar"
$ws.Range("L566").Value = "Maduro"
$ws.Range("M566").Value = 100
$ws.Range("N566").Value = 21000
$ws.Range("O566").Value = 21000
$ws.Range("P566").Value = 21000
$ws.Range("Q566").Value = "$/caja 20 kilos"
$ws.Range("R566").Value = "Ecuador"
$ws.Range("S566").Value = 1050
$ws.Range("T566").Value = 20

$ws.Range("A567").Value = 11
$ws.Range("B567").Value = "Vega Monumental Concepción"
$ws.Range("C567").Value = "Bíobío"
$ws.Range("D567").Value = 44964
$ws.Range("E567").Value = 8
$ws.Range("F567").Value = "Fruta"
$ws.Range("G567").Value = 100108
$ws.Range("H567").Value = "Tropicales y subtropicales"
$ws.Range("I567").Value = 100108006
$ws.Range("J567").Value = "Plátano"
$ws.Range("K567").Value = "Sin especificar"
$ws.Range("L567").Value = "Pintón"
$ws.Range("M567").Value = 300
$ws.Range("N567").Value = 22000
$ws.Range("O567").Value = 22000
$ws.Range("P567").Value = 22000
$ws.Range("Q567").Value = "$/caja 20 kilos"
$ws.Range("R567").Value = "Ecuador"
$ws.Range("S567").Value = 1100
$ws.Range("T567").Value = 20

$ws.Range("A568").Value = 11
$ws.Range("B568").Value = "Vega Monumental Concepción"
$ws.Range("C568").Value = "Bíobío"
$ws.Range("D568").Value = 44964
$ws.Range("E568").Value = 8
$ws.Range("F568").Value = "Fruta"
$ws.Range("G568").Value = 100108
$ws.Range("H568").Value = "Tropicales y subtropicales"
$ws.Range("I568").Value = 100108006
$ws.Range("J568").Value = "Plátano"
$ws.Range("K568").Value = "Sin especificar"
$ws.Range("L568").Value = "Primera Pintón"
$ws.Range("M568").Value = 300
$ws.Range("N568").Value = 24000
$ws.Range("O568").Value = 24000
$ws.Range("P568").Value = 24000
$ws.Range("Q568").Value = "$/caja 20 kilos"
$ws.Range("R568").Value = "Ecuador"
$ws.Range("S568").Value = 1200
$ws.Range("T568").Value = 20
